$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.842606
$ws.Range("H2").Value = 2.527818
$ws.Range("I2").Value = 0.0108780433452729
$ws.Range("J2").Value = 0.0108780433452729
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.438062
$ws.Range("N2").Value = 1.314186
$ws.Range("O2").Value = 0.6074000808827777
$ws.Range("P2").Value = 0.6074000808827777
$ws.Range("Q2").Value = 0.369113669572
$ws.Range("R2").Value = 3.322023026148
$ws.Range("S2").Value = 0.006607324407765123
$ws.Range("T2").Value = 0.006607324407765123

# Row 3 updates
$ws.Range("A3").Value = "ECs"
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.842606
$ws.Range("H3").Value = 2.527818
$ws.Range("I3").Value = 0.0108780433452729
$ws.Range("J3").Value = 0.0108780433452729
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.2831463333333333
$ws.Range("N3").Value = 0.849439
$ws.Range("O3").Value = 0.3925999191172223
$ws.Range("P3").Value = 0.3925999191172223
$ws.Range("Q3").Value = 0.2385807993446667
$ws.Range("R3").Value = 2.147227194102
$ws.Range("S3").Value = 0.004270718937507779
$ws.Range("T3").Value = 0.004270718937507779

# Row 4 updates
$ws.Range("A4").Value = "FAPs"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 69.05064766666666
$ws.Range("H4").Value = 207.151943
$ws.Range("I4").Value = 0.8914438519749055
$ws.Range("J4").Value = 0.8914438519749054
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.438062
$ws.Range("N4").Value = 1.314186
$ws.Range("O4").Value = 0.6074000808827777
$ws.Range("P4").Value = 0.6074000808827777
$ws.Range("Q4").Value = 30.24846481815533
$ws.Range("R4").Value = 272.236183363398
$ws.Range("S4").Value = 0.5414630677920126
$ws.Range("T4").Value = 0.5414630677920125

# Row 5 updates
$ws.Range("A5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 69.05064766666666
$ws.Range("H5").Value = 207.151943
$ws.Range("I5").Value = 0.8914438519749055
$ws.Range("J5").Value = 0.8914438519749054
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.2831463333333333
$ws.Range("N5").Value = 0.849439
$ws.Range("O5").Value = 0.3925999191172223
$ws.Range("P5").Value = 0.3925999191172223
$ws.Range("Q5").Value = 19.55143770110855
$ws.Range("R5").Value = 175.962939309977
$ws.Range("S5").Value = 0.349980784182893
$ws.Range("T5").Value = 0.3499807841828929

# Row 6 new
$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Ccl11"
$ws.Range("C6").Value = "Ccr3"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.558577
$ws.Range("H6").Value = 4.675731
$ws.Range("I6").Value = 0.02012122885778811
$ws.Range("J6").Value = 0.02012122885778811
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.438062
$ws.Range("N6").Value = 1.314186
$ws.Range("O6").Value = 0.6074000808827777
$ws.Range("P6").Value = 0.6074000808827777
$ws.Range("Q6").Value = 0.6827533577739999
$ws.Range("R6").Value = 6.144780219966
$ws.Range("S6").Value = 0.01222163603568138
$ws.Range("T6").Value = 0.01222163603568138

# Row 7 new
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Ccl11"
$ws.Range("C7").Value = "Ccr3"
$ws.Range("D7").Value = "M2"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.558577
$ws.Range("H7").Value = 4.675731
$ws.Range("I7").Value = 0.02012122885778811
$ws.Range("J7").Value = 0.02012122885778811
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.2831463333333333
$ws.Range("N7").Value = 0.849439
$ws.Range("O7").Value = 0.3925999191172223
$ws.Range("P7").Value = 0.3925999191172223
$ws.Range("Q7").Value = 0.4413053627676666
$ws.Range("R7").Value = 3.971748264909
$ws.Range("S7").Value = 0.007899592822106729
$ws.Range("T7").Value = 0.007899592822106729

# Row 8 new
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Ccl11"
$ws.Range("C8").Value = "Ccr3"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 6.007504
$ws.Range("H8").Value = 18.022512
$ws.Range("I8").Value = 0.07755687582203348
$ws.Range("J8").Value = 0.07755687582203348
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.438062
$ws.Range("N8").Value = 1.314186
$ws.Range("O8").Value = 0.6074000808827777
$ws.Range("P8").Value = 0.6074000808827777
$ws.Range("Q8").Value = 2.631659217248
$ws.Range("R8").Value = 23.684932955232
$ws.Range("S8").Value = 0.04710805264731869
$ws.Range("T8").Value = 0.04710805264731869

# Row 9 new
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Ccl11"
$ws.Range("C9").Value = "Ccr3"
$ws.Range("D9").Value = "M2"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 6.007504
$ws.Range("H9").Value = 18.022512
$ws.Range("I9").Value = 0.07755687582203348
$ws.Range("J9").Value = 0.07755687582203348
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.2831463333333333
$ws.Range("N9").Value = 0.849439
$ws.Range("O9").Value = 0.3925999191172223
$ws.Range("P9").Value = 0.3925999191172223
$ws.Range("Q9").Value = 1.701002730085333
$ws.Range("R9").Value = 15.309024570768
$ws.Range("S9").Value = 0.03044882317471479
$ws.Range("T9").Value = 0.03044882317471479
